$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for I0 (I1) and IF (J1), matching the formatting used by
# the existing header cells (B1:H1): bold font, thin border, centered
# horizontally, top-aligned vertically.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Row => (I value, J value) pairs taken from the diff
$data = @(
    @(2, 8, 8),
    @(3, 8, 9),
    @(4, 7, 8),
    @(5, 9, 9),
    @(6, 7, 8),
    @(7, 7, 8),
    @(8, 7, 8),
    @(9, 8, 8),
    @(10, 8, 8),
    @(11, 8, 8),
    @(12, 7, 8),
    @(13, 9, 9),
    @(14, 7, 7),
    @(15, 9, 9),
    @(16, 9, 10),
    @(17, 8, 8),
    @(18, 8, 8),
    @(19, 7, 8),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 6, 7),
    @(23, 9, 9),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 8, 8),
    @(27, 8, 9),
    @(28, 7, 7),
    @(29, 9, 9),
    @(30, 9, 9),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 7, 7),
    @(34, 7, 7),
    @(35, 7, 8),
    @(36, 5, 6),
    @(37, 8, 8),
    @(38, 7, 7),
    @(39, 5, 5),
    @(40, 7, 7),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 10, 10),
    @(44, 8, 8),
    @(45, 7, 7),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 8, 8),
    @(49, 4, 4),
    @(50, 7, 7),
    @(51, 6, 6),
    @(52, 6, 6),
    @(53, 7, 7),
    @(54, 8, 8),
    @(55, 8, 8),
    @(56, 7, 7),
    @(57, 8, 8),
    @(58, 8, 8),
    @(59, 6, 6),
    @(60, 6, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
